# lesson 498 - Monday
#
# 1. Remove the stray _GoBack bookmark that sat after "surgical" in the
#    "US army has conducted ... surgical" line.
# 2. Split three ellipsis runs by inserting a word right after the first
#    ellipsis character:
#      "put up a .......... fight"                 -> "... rough .......... fight"
#      "finally lost the ................."         -> "... battle ................"
#      "power in the country has been ..........."  -> "... disrupted .........."
# 3. Re-create the _GoBack bookmark right after the newly-typed "disrupted".
#
# Track Changes is toggled on for the text-splitting inserts and then
# accepted immediately; the engine only keeps freshly-inserted text as a
# separate run (matching the source diff) when it first arrives wrapped in
# a <w:ins>, so this is how a genuine split into sibling runs with
# identical rPr survives the save.

$d = $word.ActiveDocument

# --- 1. drop the old _GoBack bookmark -------------------------------------
$old = $d.Bookmarks.Item("_GoBack")
$old.Delete()

$d.TrackRevisions = $true

# --- 2a. "put up a .......... fight" -> insert "rough" after first dot ----
$r = $d.Content
$r.Find.Execute("put up a ……………………….. fight against", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$prefixLen = "put up a ".Length
$splitPoint = $r.Start + $prefixLen + 1
$ip = $d.Range($splitPoint, $splitPoint)
$ip.InsertBefore("rough")

# --- 2b. "finally lost the ................." -> insert "battle" --------
$r = $d.Content
$r.Find.Execute("finally lost the ……………………………….", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$prefixLen = "finally lost the ".Length
$splitPoint = $r.Start + $prefixLen + 1
$ip = $d.Range($splitPoint, $splitPoint)
$ip.InsertBefore("battle")

# --- 2c. "has been ...........by rebellions" -> insert "disrupted" -------
$r = $d.Content
$r.Find.Execute("has been ………………………..by rebellions", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$prefixLen = "has been ".Length
$ellipsisStart = $r.Start + $prefixLen
$splitPoint = $ellipsisStart + 1
$ip = $d.Range($splitPoint, $splitPoint)
$ip.InsertBefore("disrupted")

$d.TrackRevisions = $false
$d.AcceptAllRevisions()

# --- 3. re-create _GoBack right after "disrupted" -------------------------
$bmPos = $splitPoint + "disrupted".Length
$bmRange = $d.Range($bmPos, $bmPos)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null
